# shop_1030am.xlsx edit — "age, address and decode time"
#
# The underlying change (per the OOXML diff) is:
#   - F2:F6 / I2:I6 (opening/closing time columns) are converted from real
#     Excel time serial numbers (formatted h:mm AM/PM) into literal TEXT
#     that spells out the time ("10:00:00" / "10:30:00", most of them with
#     a trailing space), entered with a leading apostrophe (quote-prefix)
#     so Excel keeps them as text instead of re-parsing them as numbers.
#     This flips the cell style's quotePrefix flag, exactly like typing
#     '10:00:00  into the cell in the real app.
#   - Rows 2-6 drop their explicit (wrapped-text) row height and go back
#     to auto/default height now that the time values are short text.
#   - The selection/scroll position of the sheet view also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F column ("opening_time") -> quote-prefixed text "10:00:00" ------
# Row 6 was retyped without the trailing space; the rest keep it.
$ws.Range("F2").Value = "'10:00:00 "
$ws.Range("F3").Value = "'10:00:00 "
$ws.Range("F4").Value = "'10:00:00 "
$ws.Range("F5").Value = "'10:00:00 "
$ws.Range("F6").Value = "'10:00:00"

# --- I column ("closing time") -> quote-prefixed text "10:30:00" ------
# Row 4 was retyped without the trailing space; the rest keep it.
$ws.Range("I2").Value = "'10:30:00 "
$ws.Range("I3").Value = "'10:30:00 "
$ws.Range("I4").Value = "'10:30:00"
$ws.Range("I5").Value = "'10:30:00 "
$ws.Range("I6").Value = "'10:30:00 "

# --- Rows 2-6: back to default/auto row height (no more forced 28.8) --
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()

# --- Sheet-level default column width (19.109375 chars) ---------------
# Best-effort: mirrors sheetFormatPr/@defaultColWidth from the target file.
$ws.StandardWidth = 19.109375

# --- View state: scrolled so column C is left-most, I19 selected ------
# Best-effort: matches sheetView/@topLeftCell + <selection> in the target.
$ws.Range("I19").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
